# Update FTSE 100 ticker symbols (data refresh 2024-05-26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "AV"
$ws.Range("A13").Value = "BA"
$ws.Range("A18").Value = "BP"
$ws.Range("A20").Value = "BT-A"
$ws.Range("A52").Value = "JD"
$ws.Range("A76").Value = "RR"
$ws.Range("A93").Value = "TW"
$ws.Range("A96").Value = "UU"
